# Auto-generated edit script: update TPM-derived NATMI LR-pair stats for Pf4-Cxcr3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2859746666666667
$ws.Range("H2").Value = 0.857924
$ws.Range("I2").Value = 0.001281503322100235
$ws.Range("J2").Value = 0.001281503322100236
$ws.Range("O2").Value = 0.01611173663836548
$ws.Range("P2").Value = 0.01611173663836548
$ws.Range("Q2").Value = 0.01821687224133333
$ws.Range("R2").Value = 0.163951850172
$ws.Range("S2").Value = 2064724402686944.0/100000000000000000000.0
$ws.Range("T2").Value = 2064724402686944.0/100000000000000000000.0
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2859746666666667
$ws.Range("H3").Value = 0.857924
$ws.Range("I3").Value = 0.001281503322100235
$ws.Range("J3").Value = 0.001281503322100236
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.701496333333333
$ws.Range("N3").Value = 8.104489000000001
$ws.Range("O3").Value = 0.68328279700753
$ws.Range("P3").Value = 0.68328279700753
$ws.Range("Q3").Value = 0.7725595134262222
$ws.Range("R3").Value = 6.953035620836001
$ws.Range("S3").Value = 0.0008756291742990905
$ws.Range("T3").Value = 0.0008756291742990906
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.2859746666666667
$ws.Range("H4").Value = 0.857924
$ws.Range("I4").Value = 0.001281503322100235
$ws.Range("J4").Value = 0.001281503322100236
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.188504333333333
$ws.Range("N4").Value = 3.565513
$ws.Range("O4").Value = 0.3006054663541045
$ws.Range("P4").Value = 0.3006054663541044
$ws.Range("Q4").Value = 0.3398821305568889
$ws.Range("R4").Value = 3.058939175012
$ws.Range("S4").Value = 0.0003852269037742754
$ws.Range("T4").Value = 0.0003852269037742754
$ws.Range("G5").Value = 0.8039883333333333
$ws.Range("I5").Value = 0.003602814655248594
$ws.Range("J5").Value = 0.003602814655248594
$ws.Range("O5").Value = 0.01611173663836548
$ws.Range("P5").Value = 0.01611173663836548
$ws.Range("Q5").Value = 0.05121486082166666
$ws.Range("S5").Value = 5804760088220887.0/100000000000000000000.0
$ws.Range("T5").Value = 5804760088220885.0/100000000000000000000.0
$ws.Range("G6").Value = 0.8039883333333333
$ws.Range("I6").Value = 0.003602814655248594
$ws.Range("J6").Value = 0.003602814655248594
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.701496333333333
$ws.Range("N6").Value = 8.104489000000001
$ws.Range("O6").Value = 0.68328279700753
$ws.Range("P6").Value = 0.68328279700753
$ws.Range("Q6").Value = 2.171971534542778
$ws.Range("R6").Value = 19.547743810885
$ws.Range("S6").Value = 0.002461741274737979
$ws.Range("T6").Value = 0.002461741274737979
$ws.Range("G7").Value = 0.8039883333333333
$ws.Range("I7").Value = 0.003602814655248594
$ws.Range("J7").Value = 0.003602814655248594
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.188504333333333
$ws.Range("N7").Value = 3.565513
$ws.Range("O7").Value = 0.3006054663541045
$ws.Range("P7").Value = 0.3006054663541044
$ws.Range("Q7").Value = 0.9555436181161111
$ws.Range("R7").Value = 8.599892563045
$ws.Range("S7").Value = 0.001083025779628406
$ws.Range("T7").Value = 0.001083025779628406
$ws.Range("G8").Value = 54.865331
$ws.Range("H8").Value = 164.595993
$ws.Range("I8").Value = 0.2458613022061244
$ws.Range("J8").Value = 0.2458613022061244
$ws.Range("O8").Value = 0.01611173663836548
$ws.Range("P8").Value = 0.01611173663836548
$ws.Range("Q8").Value = 3.494976450030999
$ws.Range("R8").Value = 31.454788050279
$ws.Range("S8").Value = 0.003961252550710662
$ws.Range("T8").Value = 0.003961252550710661
$ws.Range("G9").Value = 54.865331
$ws.Range("H9").Value = 164.595993
$ws.Range("I9").Value = 0.2458613022061244
$ws.Range("J9").Value = 0.2458613022061244
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.701496333333333
$ws.Range("N9").Value = 8.104489000000001
$ws.Range("O9").Value = 0.68328279700753
$ws.Range("P9").Value = 0.68328279700753
$ws.Range("Q9").Value = 148.2184905236197
$ws.Range("R9").Value = 1333.966414712577
$ws.Range("S9").Value = 0.1679927982473143
$ws.Range("T9").Value = 0.1679927982473143
$ws.Range("G10").Value = 54.865331
$ws.Range("H10").Value = 164.595993
$ws.Range("I10").Value = 0.2458613022061244
$ws.Range("J10").Value = 0.2458613022061244
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.188504333333333
$ws.Range("N10").Value = 3.565513
$ws.Range("O10").Value = 0.3006054663541045
$ws.Range("P10").Value = 0.3006054663541044
$ws.Range("Q10").Value = 65.20768364326767
$ws.Range("R10").Value = 586.869152789409
$ws.Range("S10").Value = 0.07390725140809945
$ws.Range("T10").Value = 0.07390725140809944
$ws.Range("G11").Value = 0.04972033333333333
$ws.Range("H11").Value = 0.149161
$ws.Range("I11").Value = 0.0002228056529806757
$ws.Range("J11").Value = 0.0002228056529806757
$ws.Range("O11").Value = 0.01611173663836548
$ws.Range("P11").Value = 0.01611173663836548
$ws.Range("Q11").Value = 0.003167234953666666
$ws.Range("R11").Value = 0.028505114583
$ws.Range("S11").Value = 3589786002363698.0/1000000000000000000000.0
$ws.Range("T11").Value = 3589786002363697.0/1000000000000000000000.0
$ws.Range("G12").Value = 0.04972033333333333
$ws.Range("H12").Value = 0.149161
$ws.Range("I12").Value = 0.0002228056529806757
$ws.Range("J12").Value = 0.0002228056529806757
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 2.701496333333333
$ws.Range("N12").Value = 8.104489000000001
$ws.Range("O12").Value = 0.68328279700753
$ws.Range("P12").Value = 0.68328279700753
$ws.Range("Q12").Value = 0.1343192981921111
$ws.Range("R12").Value = 1.208873683729
$ws.Range("S12").Value = 0.0001522392697577252
$ws.Range("T12").Value = 0.0001522392697577252
$ws.Range("G13").Value = 0.04972033333333333
$ws.Range("H13").Value = 0.149161
$ws.Range("I13").Value = 0.0002228056529806757
$ws.Range("J13").Value = 0.0002228056529806757
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.188504333333333
$ws.Range("N13").Value = 3.565513
$ws.Range("O13").Value = 0.3006054663541045
$ws.Range("P13").Value = 0.3006054663541044
$ws.Range("Q13").Value = 0.05909283162144444
$ws.Range("R13").Value = 0.531835484593
$ws.Range("S13").Value = 6697659722058679.0/100000000000000000000.0
$ws.Range("T13").Value = 6697659722058678.0/100000000000000000000.0
$ws.Range("G14").Value = 167.1506043333333
$ws.Range("H14").Value = 501.451813
$ws.Range("I14").Value = 0.7490315741635462
$ws.Range("J14").Value = 0.749031574163546
$ws.Range("O14").Value = 0.01611173663836548
$ws.Range("P14").Value = 0.01611173663836548
$ws.Range("Q14").Value = 10.64766064663767
$ws.Range("R14").Value = 95.828945819739
$ws.Range("S14").Value = 0.01206819945674338
$ws.Range("T14").Value = 0.01206819945674337
$ws.Range("G15").Value = 167.1506043333333
$ws.Range("H15").Value = 501.451813
$ws.Range("I15").Value = 0.7490315741635462
$ws.Range("J15").Value = 0.749031574163546
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.701496333333333
$ws.Range("N15").Value = 8.104489000000001
$ws.Range("O15").Value = 0.68328279700753
$ws.Range("P15").Value = 0.68328279700753
$ws.Range("Q15").Value = 451.5567447209509
$ws.Range("R15").Value = 4064.010702488557
$ws.Range("S15").Value = 0.511800389041421
$ws.Range("T15").Value = 0.5118003890414209
$ws.Range("G16").Value = 167.1506043333333
$ws.Range("H16").Value = 501.451813
$ws.Range("I16").Value = 0.7490315741635462
$ws.Range("J16").Value = 0.749031574163546
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.188504333333333
$ws.Range("N16").Value = 3.565513
$ws.Range("O16").Value = 0.3006054663541045
$ws.Range("P16").Value = 0.3006054663541044
$ws.Range("Q16").Value = 198.6592175694521
$ws.Range("R16").Value = 1787.932958125069
$ws.Range("S16").Value = 0.2251629856653818
$ws.Range("T16").Value = 0.2251629856653817
